# Remove the "Student code: 79101" row from the student-info table.
$d = $word.ActiveDocument

foreach ($t in $d.Tables) {
    for ($i = $t.Rows.Count; $i -ge 1; $i--) {
        $row = $t.Rows.Item($i)
        $label = $row.Cells.Item(1).Range.Text
        if ($label -like "*Student code:*") {
            $row.Delete()
        }
    }
}
